$d = $word.ActiveDocument

# --- Edit 1: Paragraph 1 text change ---
# Original: "This is a Microsoft word document."
# New: "This is a Microsoft word document.  " + red "(This is a change – Version for branch alternate)"
#       split across 3 runs matching the original authored edit.

$p1 = $d.Paragraphs.Item(1)
$p1End = $p1.Range.End - 1  # position just before the paragraph mark

# 1a. Append two trailing spaces (keeps default/no run formatting)
$ip = $d.Range($p1End, $p1End)
$ip.InsertAfter("  ")

# 1b. Append first chunk of red text: "(This is a change – Ve"
$pos = $d.Paragraphs.Item(1).Range.End - 1
$startA = $pos
$ipA = $d.Range($pos, $pos)
$ipA.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$endA = $d.Paragraphs.Item(1).Range.End - 1
$rA = $d.Range($startA, $endA)
$rA.Font.Color = 192

# 1c. Append second chunk of red text: "rsion for branch alternate"
$startB = $d.Paragraphs.Item(1).Range.End - 1
$ipB = $d.Range($startB, $startB)
$ipB.InsertAfter("rsion for branch alternate")
$endB = $d.Paragraphs.Item(1).Range.End - 1
$rB = $d.Range($startB, $endB)
$rB.Font.Color = 192

# 1d. Append closing paren: ")"
$startC = $d.Paragraphs.Item(1).Range.End - 1
$ipC = $d.Range($startC, $startC)
$ipC.InsertAfter(")")
$endC = $d.Paragraphs.Item(1).Range.End - 1
$rC = $d.Range($startC, $endC)
$rC.Font.Color = 192

# --- Edit 2: add a new empty shaded paragraph after the last paragraph ---
$endPoint = $d.Content.End - 1
$ip2 = $d.Range($endPoint, $endPoint)

$xmlFrag = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$ip2.InsertXML($xmlFrag)
